$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 79 (shifts existing rows 79-178 down to 80-179,
# carrying formatting - e.g. the date style on column D - from the row above).
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new weekly record.
$ws.Range("A79").Value = 10
$ws.Range("B79").Value = "Vega Modelo de Temuco"
$ws.Range("C79").Value = "La Araucanía"
$ws.Range("D79").Value = 44467
$ws.Range("E79").Value = 9
$ws.Range("F79").Value = 100112017
$ws.Range("G79").Value = "Apio"
$ws.Range("H79").Value = "Americana (o)"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 80
$ws.Range("K79").Value = 10000
$ws.Range("L79").Value = 10000
$ws.Range("M79").Value = 10000
$ws.Range("N79").Value = "$/docena de matas"
$ws.Range("O79").Value = "Provincia del Elquí"
$ws.Range("P79").Value = 1667
$ws.Range("Q79").Value = 6
$ws.Range("R79").Value = "Hortaliza"
